$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing data block first so stale shared-string references
# do not linger, then rebuild it from scratch.
$ws.Range("A2:T13").ClearContents()

# Step 1: write sending/target cluster + ligand/receptor symbol columns (A-D)
# in column-major order so that shared-string indices come out in the same
# order as the target workbook (M2 ends up right after FAPs).
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(17,1).Value = "sCs"

$ws.Cells.Item(2,2).Value = "Vtn"
$ws.Cells.Item(3,2).Value = "Vtn"
$ws.Cells.Item(4,2).Value = "Vtn"
$ws.Cells.Item(5,2).Value = "Vtn"
$ws.Cells.Item(6,2).Value = "Vtn"
$ws.Cells.Item(7,2).Value = "Vtn"
$ws.Cells.Item(8,2).Value = "Vtn"
$ws.Cells.Item(9,2).Value = "Vtn"
$ws.Cells.Item(10,2).Value = "Vtn"
$ws.Cells.Item(11,2).Value = "Vtn"
$ws.Cells.Item(12,2).Value = "Vtn"
$ws.Cells.Item(13,2).Value = "Vtn"
$ws.Cells.Item(14,2).Value = "Vtn"
$ws.Cells.Item(15,2).Value = "Vtn"
$ws.Cells.Item(16,2).Value = "Vtn"
$ws.Cells.Item(17,2).Value = "Vtn"

$ws.Cells.Item(2,3).Value = "Itgav"
$ws.Cells.Item(3,3).Value = "Itgav"
$ws.Cells.Item(4,3).Value = "Itgav"
$ws.Cells.Item(5,3).Value = "Itgav"
$ws.Cells.Item(6,3).Value = "Itgav"
$ws.Cells.Item(7,3).Value = "Itgav"
$ws.Cells.Item(8,3).Value = "Itgav"
$ws.Cells.Item(9,3).Value = "Itgav"
$ws.Cells.Item(10,3).Value = "Itgav"
$ws.Cells.Item(11,3).Value = "Itgav"
$ws.Cells.Item(12,3).Value = "Itgav"
$ws.Cells.Item(13,3).Value = "Itgav"
$ws.Cells.Item(14,3).Value = "Itgav"
$ws.Cells.Item(15,3).Value = "Itgav"
$ws.Cells.Item(16,3).Value = "Itgav"
$ws.Cells.Item(17,3).Value = "Itgav"

$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(17,4).Value = "sCs"

# Step 2: write the numeric columns (E-T) for every row
# row 2: ECs -> ECs
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 12.75206033333333
$ws.Cells.Item(2,8).Value = 38.256181
$ws.Cells.Item(2,9).Value = 0.1573122343381959
$ws.Cells.Item(2,10).Value = 0.157312234338196
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 13.441269
$ws.Cells.Item(2,14).Value = 40.323807
$ws.Cells.Item(2,15).Value = 0.0897308213348123
$ws.Cells.Item(2,16).Value = 0.08973082133481232
$ws.Cells.Item(2,17).Value = 171.403873244563
$ws.Cells.Item(2,18).Value = 1542.634859201067
$ws.Cells.Item(2,19).Value = 0.01411575599318078
$ws.Cells.Item(2,20).Value = 0.01411575599318079

# row 3: ECs -> FAPs
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 12.75206033333333
$ws.Cells.Item(3,8).Value = 38.256181
$ws.Cells.Item(3,9).Value = 0.1573122343381959
$ws.Cells.Item(3,10).Value = 0.157312234338196
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 54.711535
$ws.Cells.Item(3,14).Value = 164.134605
$ws.Cells.Item(3,15).Value = 0.3652416280068742
$ws.Cells.Item(3,16).Value = 0.3652416280068742
$ws.Cells.Item(3,17).Value = 697.6847952492783
$ws.Cells.Item(3,18).Value = 6279.163157243504
$ws.Cells.Item(3,19).Value = 0.05745697657508157
$ws.Cells.Item(3,20).Value = 0.05745697657508159

# row 4: ECs -> M2
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 12.75206033333333
$ws.Cells.Item(4,8).Value = 38.256181
$ws.Cells.Item(4,9).Value = 0.1573122343381959
$ws.Cells.Item(4,10).Value = 0.157312234338196
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 63.67711
$ws.Cells.Item(4,14).Value = 191.03133
$ws.Cells.Item(4,15).Value = 0.4250937452800914
$ws.Cells.Item(4,16).Value = 0.4250937452800915
$ws.Cells.Item(4,17).Value = 812.0143485723033
$ws.Cells.Item(4,18).Value = 7308.12913715073
$ws.Cells.Item(4,19).Value = 0.06687244687320311
$ws.Cells.Item(4,20).Value = 0.06687244687320314

# row 5: ECs -> sCs
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 12.75206033333333
$ws.Cells.Item(5,8).Value = 38.256181
$ws.Cells.Item(5,9).Value = 0.1573122343381959
$ws.Cells.Item(5,10).Value = 0.157312234338196
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 17.96553866666667
$ws.Cells.Item(5,14).Value = 53.896616
$ws.Cells.Item(5,15).Value = 0.119933805378222
$ws.Cells.Item(5,16).Value = 0.119933805378222
$ws.Cells.Item(5,17).Value = 229.0976329981662
$ws.Cells.Item(5,18).Value = 2061.878696983496
$ws.Cells.Item(5,19).Value = 0.01886705489673044
$ws.Cells.Item(5,20).Value = 0.01886705489673045

# row 6: FAPs -> ECs
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 20.35396833333334
$ws.Cells.Item(6,8).Value = 61.06190500000001
$ws.Cells.Item(6,9).Value = 0.2510910513649196
$ws.Cells.Item(6,10).Value = 0.2510910513649196
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 13.441269
$ws.Cells.Item(6,14).Value = 40.323807
$ws.Cells.Item(6,15).Value = 0.0897308213348123
$ws.Cells.Item(6,16).Value = 0.08973082133481232
$ws.Cells.Item(6,17).Value = 273.583163585815
$ws.Cells.Item(6,18).Value = 2462.248472272336
$ws.Cells.Item(6,19).Value = 0.02253060626879578
$ws.Cells.Item(6,20).Value = 0.02253060626879578

# row 7: FAPs -> FAPs
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 20.35396833333334
$ws.Cells.Item(7,8).Value = 61.06190500000001
$ws.Cells.Item(7,9).Value = 0.2510910513649196
$ws.Cells.Item(7,10).Value = 0.2510910513649196
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 54.711535
$ws.Cells.Item(7,14).Value = 164.134605
$ws.Cells.Item(7,15).Value = 0.3652416280068742
$ws.Cells.Item(7,16).Value = 0.3652416280068742
$ws.Cells.Item(7,17).Value = 1113.596850858059
$ws.Cells.Item(7,18).Value = 10022.37165772253
$ws.Cells.Item(7,19).Value = 0.0917089043784809
$ws.Cells.Item(7,20).Value = 0.09170890437848091

# row 8: FAPs -> M2
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 20.35396833333334
$ws.Cells.Item(8,8).Value = 61.06190500000001
$ws.Cells.Item(8,9).Value = 0.2510910513649196
$ws.Cells.Item(8,10).Value = 0.2510910513649196
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 63.67711
$ws.Cells.Item(8,14).Value = 191.03133
$ws.Cells.Item(8,15).Value = 0.4250937452800914
$ws.Cells.Item(8,16).Value = 0.4250937452800915
$ws.Cells.Item(8,17).Value = 1296.081880498184
$ws.Cells.Item(8,18).Value = 11664.73692448365
$ws.Cells.Item(8,19).Value = 0.1067372354310295
$ws.Cells.Item(8,20).Value = 0.1067372354310295

# row 9: FAPs -> sCs
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 20.35396833333334
$ws.Cells.Item(9,8).Value = 61.06190500000001
$ws.Cells.Item(9,9).Value = 0.2510910513649196
$ws.Cells.Item(9,10).Value = 0.2510910513649196
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 17.96553866666667
$ws.Cells.Item(9,14).Value = 53.896616
$ws.Cells.Item(9,15).Value = 0.119933805378222
$ws.Cells.Item(9,16).Value = 0.119933805378222
$ws.Cells.Item(9,17).Value = 365.670005112609
$ws.Cells.Item(9,18).Value = 3291.030046013481
$ws.Cells.Item(9,19).Value = 0.03011430528661341
$ws.Cells.Item(9,20).Value = 0.03011430528661341

# row 10: M2 -> ECs
$ws.Cells.Item(10,5).Value = 1
$ws.Cells.Item(10,6).Value = 0.3333333333333333
$ws.Cells.Item(10,7).Value = 0.004706
$ws.Cells.Item(10,8).Value = 0.014118
$ws.Cells.Item(10,9).Value = 0.00005805425597465284
$ws.Cells.Item(10,10).Value = 0.00005805425597465285
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 13.441269
$ws.Cells.Item(10,14).Value = 40.323807
$ws.Cells.Item(10,15).Value = 0.0897308213348123
$ws.Cells.Item(10,16).Value = 0.08973082133481232
$ws.Cells.Item(10,17).Value = 0.063254611914
$ws.Cells.Item(10,18).Value = 0.5692915072260001
$ws.Cells.Item(10,19).Value = 0.000005209256070587034
$ws.Cells.Item(10,20).Value = 0.000005209256070587036

# row 11: M2 -> FAPs
$ws.Cells.Item(11,5).Value = 1
$ws.Cells.Item(11,6).Value = 0.3333333333333333
$ws.Cells.Item(11,7).Value = 0.004706
$ws.Cells.Item(11,8).Value = 0.014118
$ws.Cells.Item(11,9).Value = 0.00005805425597465284
$ws.Cells.Item(11,10).Value = 0.00005805425597465285
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 54.711535
$ws.Cells.Item(11,14).Value = 164.134605
$ws.Cells.Item(11,15).Value = 0.3652416280068742
$ws.Cells.Item(11,16).Value = 0.3652416280068742
$ws.Cells.Item(11,17).Value = 0.25747248371
$ws.Cells.Item(11,18).Value = 2.31725235339
$ws.Cells.Item(11,19).Value = 0.00002120383096491
$ws.Cells.Item(11,20).Value = 0.00002120383096491001

# row 12: M2 -> M2
$ws.Cells.Item(12,5).Value = 1
$ws.Cells.Item(12,6).Value = 0.3333333333333333
$ws.Cells.Item(12,7).Value = 0.004706
$ws.Cells.Item(12,8).Value = 0.014118
$ws.Cells.Item(12,9).Value = 0.00005805425597465284
$ws.Cells.Item(12,10).Value = 0.00005805425597465285
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 63.67711
$ws.Cells.Item(12,14).Value = 191.03133
$ws.Cells.Item(12,15).Value = 0.4250937452800914
$ws.Cells.Item(12,16).Value = 0.4250937452800915
$ws.Cells.Item(12,17).Value = 0.29966447966
$ws.Cells.Item(12,18).Value = 2.69698031694
$ws.Cells.Item(12,19).Value = 0.0000246785011017143
$ws.Cells.Item(12,20).Value = 0.00002467850110171431

# row 13: M2 -> sCs
$ws.Cells.Item(13,5).Value = 1
$ws.Cells.Item(13,6).Value = 0.3333333333333333
$ws.Cells.Item(13,7).Value = 0.004706
$ws.Cells.Item(13,8).Value = 0.014118
$ws.Cells.Item(13,9).Value = 0.00005805425597465284
$ws.Cells.Item(13,10).Value = 0.00005805425597465285
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 17.96553866666667
$ws.Cells.Item(13,14).Value = 53.896616
$ws.Cells.Item(13,15).Value = 0.119933805378222
$ws.Cells.Item(13,16).Value = 0.119933805378222
$ws.Cells.Item(13,17).Value = 0.08454582496533333
$ws.Cells.Item(13,18).Value = 0.760912424688
$ws.Cells.Item(13,19).Value = 0.000006962667837441496
$ws.Cells.Item(13,20).Value = 0.000006962667837441497

# row 14: sCs -> ECs
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 47.951367
$ws.Cells.Item(14,8).Value = 143.854101
$ws.Cells.Item(14,9).Value = 0.5915386600409097
$ws.Cells.Item(14,10).Value = 0.5915386600409098
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 13.441269
$ws.Cells.Item(14,14).Value = 40.323807
$ws.Cells.Item(14,15).Value = 0.0897308213348123
$ws.Cells.Item(14,16).Value = 0.08973082133481232
$ws.Cells.Item(14,17).Value = 644.5272227647231
$ws.Cells.Item(14,18).Value = 5800.745004882508
$ws.Cells.Item(14,19).Value = 0.05307924981676514
$ws.Cells.Item(14,20).Value = 0.05307924981676516

# row 15: sCs -> FAPs
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 47.951367
$ws.Cells.Item(15,8).Value = 143.854101
$ws.Cells.Item(15,9).Value = 0.5915386600409097
$ws.Cells.Item(15,10).Value = 0.5915386600409098
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 54.711535
$ws.Cells.Item(15,14).Value = 164.134605
$ws.Cells.Item(15,15).Value = 0.3652416280068742
$ws.Cells.Item(15,16).Value = 0.3652416280068742
$ws.Cells.Item(15,17).Value = 2623.492893918345
$ws.Cells.Item(15,18).Value = 23611.43604526511
$ws.Cells.Item(15,19).Value = 0.2160545432223467
$ws.Cells.Item(15,20).Value = 0.2160545432223468

# row 16: sCs -> M2
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 47.951367
$ws.Cells.Item(16,8).Value = 143.854101
$ws.Cells.Item(16,9).Value = 0.5915386600409097
$ws.Cells.Item(16,10).Value = 0.5915386600409098
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 63.67711
$ws.Cells.Item(16,14).Value = 191.03133
$ws.Cells.Item(16,15).Value = 0.4250937452800914
$ws.Cells.Item(16,16).Value = 0.4250937452800915
$ws.Cells.Item(16,17).Value = 3053.40447110937
$ws.Cells.Item(16,18).Value = 27480.64023998433
$ws.Cells.Item(16,19).Value = 0.2514593844747571
$ws.Cells.Item(16,20).Value = 0.2514593844747571

# row 17: sCs -> sCs
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 47.951367
$ws.Cells.Item(17,8).Value = 143.854101
$ws.Cells.Item(17,9).Value = 0.5915386600409097
$ws.Cells.Item(17,10).Value = 0.5915386600409098
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 17.96553866666667
$ws.Cells.Item(17,14).Value = 53.896616
$ws.Cells.Item(17,15).Value = 0.119933805378222
$ws.Cells.Item(17,16).Value = 0.119933805378222
$ws.Cells.Item(17,17).Value = 861.4721379580241
$ws.Cells.Item(17,18).Value = 7753.249241622217
$ws.Cells.Item(17,19).Value = 0.07094548252704069
$ws.Cells.Item(17,20).Value = 0.07094548252704072

Write-Host "Done"